# Generate Report for Handoff
# Adds a new handoff entry (831053f0-66ba-499b-abad-18b05c0a1e66.md) as a
# third row on each of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileId = "831053f0-66ba-499b-abad-18b05c0a1e66"
$commitSha = "3b3395963b0abc886450566600399c35068172ac"
$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e"

function Set-EmptyTextCell($range) {
    # Writes an explicit empty-string shared-string cell (matches the
    # existing "" placeholder cells already in this workbook) instead of
    # leaving the cell completely absent, and strips the quote-prefix
    # marker that a literal apostrophe would otherwise leave behind.
    $range.Value = "'"
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$fileId.md"
$wsOverview.Range("C3").Value = ".md"
Set-EmptyTextCell $wsOverview.Range("D3")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = "2016-09-03 12:43:48"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$repoBase/$fileId.md", "", "", "e2e\$fileId.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "$fileId.454a0fdf45b1e92d9755e4ea5f0cdd4b0d92d054.zh-cn.xlf"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H3").Value = "2016-09-03 12:43:44"
Set-EmptyTextCell $wsZhCn.Range("I3")
Set-EmptyTextCell $wsZhCn.Range("J3")
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
Set-EmptyTextCell $wsZhCn.Range("L3")
$wsZhCn.Range("M3").Value = "'True"
Set-EmptyTextCell $wsZhCn.Range("N3")
$wsZhCn.Range("O3").Value = "'False"
Set-EmptyTextCell $wsZhCn.Range("P3")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$repoBase/$fileId.md", "", "", "$fileId.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "$fileId.454a0fdf45b1e92d9755e4ea5f0cdd4b0d92d054.de-de.xlf"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H3").Value = "2016-09-03 12:43:48"
Set-EmptyTextCell $wsDeDe.Range("I3")
Set-EmptyTextCell $wsDeDe.Range("J3")
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
Set-EmptyTextCell $wsDeDe.Range("L3")
$wsDeDe.Range("M3").Value = "'True"
Set-EmptyTextCell $wsDeDe.Range("N3")
$wsDeDe.Range("O3").Value = "'False"
Set-EmptyTextCell $wsDeDe.Range("P3")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$repoBase/$fileId.md", "", "", "$fileId.md") | Out-Null
